$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = 'id'
$ws.Range("B1").Value = 'name'
$ws.Range("C1").Value = 'descr'
$ws.Range("D1").Value = 'lang_code'
$ws.Range("E1").Value = 'is_active'
$ws.Range("F1").Value = 'cr_by'
$ws.Range("G1").Value = 'cr_dtimes'
$ws.Range("H1").Value = 'upd_by'
$ws.Range("I1").Value = 'upd_dtimes'
$ws.Range("J1").Value = 'is_deleted'
$ws.Range("K1").Value = 'del_dtimes'

# --- Data rows 2-6 ---
$ws.Range("A2").Value = 'login_auth'
$ws.Range("B2").Value = 'Authentification Ã  la connexion'
$ws.Range("C2").Value = 'Authentification Ã  la connexion de l''utilisateur'
$ws.Range("D2").Value = 'fra'
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = 'superadmin'
$ws.Range("G2").Value = 45079.577233923614
$ws.Range("G2").NumberFormat = "mm:ss.0"
$ws.Range("H2").Value = 'NULL'
$ws.Range("I2").Value = 'NULL'
$ws.Range("J2").Value = $false
$ws.Range("K2").Value = 'NULL'

$ws.Range("A3").Value = 'eod_auth'
$ws.Range("B3").Value = 'Authentification de Ã  la clÃ´ture'
$ws.Range("C3").Value = 'Authentification Ã  la clÃ´ture'
$ws.Range("D3").Value = 'fra'
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = 'superadmin'
$ws.Range("G3").Value = 45079.577233923614
$ws.Range("G3").NumberFormat = "mm:ss.0"
$ws.Range("H3").Value = 'NULL'
$ws.Range("I3").Value = 'NULL'
$ws.Range("J3").Value = $false
$ws.Range("K3").Value = 'NULL'

$ws.Range("A4").Value = 'packet_auth'
$ws.Range("B4").Value = 'Authentification de paquet'
$ws.Range("C4").Value = 'Authentification de paquet'
$ws.Range("D4").Value = 'fra'
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = 'superadmin'
$ws.Range("G4").Value = 45079.577233923614
$ws.Range("G4").NumberFormat = "mm:ss.0"
$ws.Range("H4").Value = 'NULL'
$ws.Range("I4").Value = 'NULL'
$ws.Range("J4").Value = $false
$ws.Range("K4").Value = 'NULL'

$ws.Range("A5").Value = 'exception_auth'
$ws.Range("B5").Value = 'Authentification dâ€™exception'
$ws.Range("C5").Value = 'Authentification dâ€™exception'
$ws.Range("D5").Value = 'fra'
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = 'superadmin'
$ws.Range("G5").Value = 45079.577233923614
$ws.Range("G5").NumberFormat = "mm:ss.0"
$ws.Range("H5").Value = 'NULL'
$ws.Range("I5").Value = 'NULL'
$ws.Range("J5").Value = $false
$ws.Range("K5").Value = 'NULL'

$ws.Range("A6").Value = 'onboard_auth'
$ws.Range("B6").Value = 'Authentification embarquÃ©e'
$ws.Range("C6").Value = 'Authentification embarquÃ©e'
$ws.Range("D6").Value = 'fra'
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = 'superadmin'
$ws.Range("G6").Value = 45079.577233923614
$ws.Range("G6").NumberFormat = "mm:ss.0"
$ws.Range("H6").Value = 'NULL'
$ws.Range("I6").Value = 'NULL'
$ws.Range("J6").Value = $false
$ws.Range("K6").Value = 'NULL'

# --- Selection matches target (B13) ---
[void]$ws.Range("B13").Select()
